$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "27.484.55"
$cell.Style = $origStyle
$cell = $ws.Range("E2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -2.28%  "
$cell.Style = $origStyle
$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.748.57"
$cell.Style = $origStyle
$cell = $ws.Range("E3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -2.41%  "
$cell.Style = $origStyle
$cell = $ws.Range("E4")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.19%  "
$cell.Style = $origStyle
$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "324.16"
$cell.Style = $origStyle
$cell = $ws.Range("E5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.27%  "
$cell.Style = $origStyle
$cell = $ws.Range("E6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.06%  "
$cell.Style = $origStyle
$cell = $ws.Range("D7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.4468"
$cell.Style = $origStyle
$cell = $ws.Range("E7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +4.41%  "
$cell.Style = $origStyle
$cell = $ws.Range("E8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.54%  "
$cell.Style = $origStyle
$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.07474"
$cell.Style = $origStyle
$cell = $ws.Range("E9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.69%  "
$cell.Style = $origStyle
$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "42.05"
$cell.Style = $origStyle
$cell = $ws.Range("E10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -5.92%  "
$cell.Style = $origStyle
$cell = $ws.Range("E11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -2.14%  "
$cell.Style = $origStyle
$cell = $ws.Range("E12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.21%  "
$cell.Style = $origStyle
$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "20.56"
$cell.Style = $origStyle
$cell = $ws.Range("E13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -4.96%  "
$cell.Style = $origStyle
$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.012"
$cell.Style = $origStyle
$cell = $ws.Range("E14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -2.24%  "
$cell.Style = $origStyle
$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.106"
$cell.Style = $origStyle
$cell = $ws.Range("E15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -3.10%  "
$cell.Style = $origStyle
$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.749.93"
$cell.Style = $origStyle
$cell = $ws.Range("E16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -3.37%  "
$cell.Style = $origStyle
$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "92.14"
$cell.Style = $origStyle
$cell = $ws.Range("E17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.23%  "
$cell.Style = $origStyle
$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.00001058"
$cell.Style = $origStyle
$cell = $ws.Range("E18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -1.20%  "
$cell.Style = $origStyle
$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.06402"
$cell.Style = $origStyle
$cell = $ws.Range("E19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.83%  "
$cell.Style = $origStyle
$cell = $ws.Range("E20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.03%  "
$cell.Style = $origStyle
$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "16.78"
$cell.Style = $origStyle
$cell = $ws.Range("E21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -2.31%  "
$cell.Style = $origStyle
$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.834"
$cell.Style = $origStyle
$cell = $ws.Range("E22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -2.34%  "
$cell.Style = $origStyle
$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "27.534.11"
$cell.Style = $origStyle
$cell = $ws.Range("E23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -2.19%  "
$cell.Style = $origStyle
$cell = $ws.Range("E24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -2.36%  "
$cell.Style = $origStyle
$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.117"
$cell.Style = $origStyle
$cell = $ws.Range("E25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -2.64%  "
$cell.Style = $origStyle
$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "161.80"
$cell.Style = $origStyle
$cell = $ws.Range("E26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.76%  "
$cell.Style = $origStyle
$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "20.40"
$cell.Style = $origStyle
$cell = $ws.Range("E27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.46%  "
$cell.Style = $origStyle
$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.951.40"
$cell.Style = $origStyle
$cell = $ws.Range("E28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -3.17%  "
$cell.Style = $origStyle
$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.071"
$cell.Style = $origStyle
$cell = $ws.Range("E29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -7.07%  "
$cell.Style = $origStyle
$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "124.45"
$cell.Style = $origStyle
$cell = $ws.Range("E30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -2.57%  "
$cell.Style = $origStyle
$cell = $ws.Range("E31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -7.87%  "
$cell.Style = $origStyle
$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.657"
$cell.Style = $origStyle
$cell = $ws.Range("E32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +3.49%  "
$cell.Style = $origStyle
$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.09012"
$cell.Style = $origStyle
$cell = $ws.Range("E33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.06%  "
$cell.Style = $origStyle
$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.486"
$cell.Style = $origStyle
$cell = $ws.Range("E34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -5.81%  "
$cell.Style = $origStyle
$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "11.97"
$cell.Style = $origStyle
$cell = $ws.Range("E35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -5.93%  "
$cell.Style = $origStyle
$cell = $ws.Range("E36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -2.48%  "
$cell.Style = $origStyle
$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.2081"
$cell.Style = $origStyle
$cell = $ws.Range("E37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -1.63%  "
$cell.Style = $origStyle
$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.6339"
$cell.Style = $origStyle
$cell = $ws.Range("E38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -2.13%  "
$cell.Style = $origStyle
$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.05974"
$cell.Style = $origStyle
$cell = $ws.Range("E39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -2.12%  "
$cell.Style = $origStyle
$cell = $ws.Range("B40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "TrustWalletToken"
$cell.Style = $origStyle
$cell = $ws.Range("C40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$cell.Style = $origStyle
$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.207"
$cell.Style = $origStyle
$cell = $ws.Range("E40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.28%  "
$cell.Style = $origStyle
$cell = $ws.Range("B41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "InternetComputer(DFINITY)"
$cell.Style = $origStyle
$cell = $ws.Range("C41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$cell.Style = $origStyle
$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.914"
$cell.Style = $origStyle
$cell = $ws.Range("E41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -3.48%  "
$cell.Style = $origStyle
$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.Style = $origStyle
$cell = $ws.Range("E42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.07%  "
$cell.Style = $origStyle
$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.386"
$cell.Style = $origStyle
$cell = $ws.Range("E43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -2.97%  "
$cell.Style = $origStyle
$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.768"
$cell.Style = $origStyle
$cell = $ws.Range("E44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -2.03%  "
$cell.Style = $origStyle
$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "13.14"
$cell.Style = $origStyle
$cell = $ws.Range("E45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -3.46%  "
$cell.Style = $origStyle
$cell = $ws.Range("B46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "PancakeSwap"
$cell.Style = $origStyle
$cell = $ws.Range("C46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$cell.Style = $origStyle
$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.703"
$cell.Style = $origStyle
$cell = $ws.Range("E46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.03%  "
$cell.Style = $origStyle
$cell = $ws.Range("B47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "Decentraland"
$cell.Style = $origStyle
$cell = $ws.Range("C47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$cell.Style = $origStyle
$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.5867"
$cell.Style = $origStyle
$cell = $ws.Range("E47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -2.35%  "
$cell.Style = $origStyle
$cell = $ws.Range("E48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -3.36%  "
$cell.Style = $origStyle
$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.943"
$cell.Style = $origStyle
$cell = $ws.Range("E49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -2.33%  "
$cell.Style = $origStyle
$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.146"
$cell.Style = $origStyle
$cell = $ws.Range("E50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.45%  "
$cell.Style = $origStyle
$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.06851"
$cell.Style = $origStyle
$cell = $ws.Range("E51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -1.67%  "
$cell.Style = $origStyle
